# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 per the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its exact textual representation (e.g. trailing
# zeros like "1.00", or thousand-dot formatted numbers like "45.382.51")
# instead of being auto-converted to a numeric value by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row = 2;  D = "45.382.51"; E = "  +4.68%  " },
    @{ Row = 3;  D = "2.369.31";  E = "  +1.83%  " },
    @{ Row = 4;  D = "1.00";      E = "  -0.63%  " },
    @{ Row = 5;  D = "310.58";    E = "  -0.55%  " },
    @{ Row = 6;  D = "108.48";    E = "  +0.11%  " },
    @{ Row = 7;  D = "0.629";     E = "  -0.23%  " },
    @{ Row = 8;  D = "1.00";      E = "  -0.32%  " },
    @{ Row = 9;  D = "0.617";     E = "  +0.63%  " },
    @{ Row = 10; D = "41.20";     E = "  +1.61%  " },
    @{ Row = 11; D = "0.0920";    E = "  +0.24%  " },
    @{ Row = 12; D = "8.51";      E = "  -0.59%  " },
    @{ Row = 13; D = "0.110";     E = "  +1.85%  " },
    @{ Row = 14; D = "0.982";     E = "  -2.09%  " },
    @{ Row = 15; D = "2.724.20";  E = "  +1.62%  " },
    @{ Row = 16; D = "15.25";     E = "  -0.90%  " },
    @{ Row = 17; D = "2.356.08";  E = "  +1.43%  " },
    @{ Row = 18; D = "45.276.52"; E = "  +4.59%  " },
    @{ Row = 19; D = "14.37";     E = "  +10.03%  " },
    @{ Row = 20; D = "7.32";      E = "  -2.47%  " },
    @{ Row = 21; D = "0.0000106"; E = "  -0.29%  " },
    @{ Row = 22; D = "73.23";     E = "  -1.05%  " },
    @{ Row = 23; D = "3.50";      E = "  +0.01%  " },
    @{ Row = 24; D = "260.66";    E = "  -2.87%  " },
    @{ Row = 25; D = "2.32";      E = "  +2.57%  " },
    @{ Row = 26; D = "1.00";      E = "  -0.30%  " },
    @{ Row = 27; D = "11.18";     E = "  +0.46%  " },
    @{ Row = 28; D = "7.36";      E = "  -2.77%  " },
    @{ Row = 29; D = "2.34";      E = "  +2.13%  " },
    @{ Row = 30; D = "0.0969";    E = "  +9.61%  " },
    @{ Row = 31; D = "22.37";     E = "  -0.96%  " },
    @{ Row = 32; D = "37.69";     E = "  -3.68%  " },
    @{ Row = 33; D = "169.15";    E = "  +1.18%  " },
    @{ Row = 34; D = "2.94";      E = "  +5.59%  " },
    @{ Row = 35; D = "0.130";     E = "  -0.49%  " },
    @{ Row = 36; D = "0.118";     E = "  +4.41%  " },
    @{ Row = 37; D = "4.79";      E = "  +1.36%  " },
    @{ Row = 38; D = "2.97";      E = "  +5.00%  " },
    @{ Row = 39; D = "3.92";      E = "  +3.78%  " },
    @{ Row = 40; D = "0.0355";    E = "  -1.51%  " },
    @{ Row = 41; D = "1.74";      E = "  +3.92%  " },
    @{ Row = 42; D = "99.50";     E = "  -5.31%  " },
    @{ Row = 43; D = "0.233";     E = "  -0.84%  " },
    @{ Row = 44; D = "69.73";     E = "  -2.70%  " },
    @{ Row = 45; D = "12.94";     E = "  -4.26%  " },
    @{ Row = 46; D = "1.00";      E = "  -0.34%  " },
    @{ Row = 47; D = "81.05";     E = "  +4.98%  " },
    @{ Row = 48; D = "112.53";    E = "  -1.07%  " },
    @{ Row = 49; D = "5.55";      E = "  +4.33%  " },
    @{ Row = 50; D = "9.25";      E = "  +3.72%  " },
    @{ Row = 51; D = "1.671.06";  E = "  +0.47%  " }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
